$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (F) values after repulling data / recalculating means
$updates = @{
    'F10' = -4
    'F12' = -3
    'F18' = -2
    'F20' = 4
    'F23' = -6
    'F27' = 0
    'F30' = -5
    'F41' = 2
    'F46' = -4
    'F47' = 4
    'F48' = -5
    'F51' = 1
    'F56' = -4
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
